{"js": "// Replace the header date and every three-digit \u00f7 one-digit math-fact\n// answer in the practice table with the new set of values.\nconst replacements = [\n  [\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"],\n  [\"857\u00f72=428, 1\", \"328\u00f77=46, 6\"],\n  [\"468\u00f74=117, 0\", \"724\u00f75=144, 4\"],\n  [\"428\u00f78=53, 4\", \"677\u00f73=225, 2\"],\n  [\"516\u00f75=103, 1\", \"309\u00f78=38, 5\"],\n  [\"892\u00f73=297, 1\", \"418\u00f72=209, 0\"],\n  [\"768\u00f77=109, 5\", \"660\u00f72=330, 0\"],\n  [\"514\u00f73=171, 1\", \"300\u00f74=75, 0\"],\n  [\"461\u00f72=230, 1\", \"694\u00f75=138, 4\"],\n  [\"936\u00f73=312, 0\", \"978\u00f78=122, 2\"],\n  [\"540\u00f79=60, 0\", \"744\u00f79=82, 6\"],\n  [\"545\u00f77=77, 6\", \"562\u00f77=80, 2\"],\n  [\"212\u00f77=30, 2\", \"571\u00f77=81, 4\"],\n  [\"466\u00f72=233, 0\", \"613\u00f75=122, 3\"],\n  [\"879\u00f74=219, 3\", \"391\u00f78=48, 7\"],\n  [\"722\u00f73=240, 2\", \"797\u00f72=398, 1\"],\n  [\"152\u00f76=25, 2\", \"351\u00f79=39, 0\"],\n  [\"484\u00f72=242, 0\", \"829\u00f76=138, 1\"],\n  [\"320\u00f76=53, 2\", \"661\u00f74=165, 1\"],\n  [\"677\u00f78=84, 5\", \"888\u00f74=222, 0\"],\n  [\"230\u00f75=46, 0\", \"373\u00f78=46, 5\"],\n  [\"409\u00f75=81, 4\", \"471\u00f79=52, 3\"],\n  [\"498\u00f75=99, 3\", \"414\u00f75=82, 4\"],\n  [\"708\u00f79=78, 6\", \"795\u00f77=113, 4\"],\n  [\"904\u00f76=150, 4\", \"939\u00f72=469, 1\"],\n  [\"102\u00f77=14, 4\", \"487\u00f78=60, 7\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the header date and every three-digit / one-digit division fact\n# answer in the practice table with the new set of values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"),\n    @(\"857\u00f72=428, 1\", \"328\u00f77=46, 6\"),\n    @(\"468\u00f74=117, 0\", \"724\u00f75=144, 4\"),\n    @(\"428\u00f78=53, 4\", \"677\u00f73=225, 2\"),\n    @(\"516\u00f75=103, 1\", \"309\u00f78=38, 5\"),\n    @(\"892\u00f73=297, 1\", \"418\u00f72=209, 0\"),\n    @(\"768\u00f77=109, 5\", \"660\u00f72=330, 0\"),\n    @(\"514\u00f73=171, 1\", \"300\u00f74=75, 0\"),\n    @(\"461\u00f72=230, 1\", \"694\u00f75=138, 4\"),\n    @(\"936\u00f73=312, 0\", \"978\u00f78=122, 2\"),\n    @(\"540\u00f79=60, 0\", \"744\u00f79=82, 6\"),\n    @(\"545\u00f77=77, 6\", \"562\u00f77=80, 2\"),\n    @(\"212\u00f77=30, 2\", \"571\u00f77=81, 4\"),\n    @(\"466\u00f72=233, 0\", \"613\u00f75=122, 3\"),\n    @(\"879\u00f74=219, 3\", \"391\u00f78=48, 7\"),\n    @(\"722\u00f73=240, 2\", \"797\u00f72=398, 1\"),\n    @(\"152\u00f76=25, 2\", \"351\u00f79=39, 0\"),\n    @(\"484\u00f72=242, 0\", \"829\u00f76=138, 1\"),\n    @(\"320\u00f76=53, 2\", \"661\u00f74=165, 1\"),\n    @(\"677\u00f78=84, 5\", \"888\u00f74=222, 0\"),\n    @(\"230\u00f75=46, 0\", \"373\u00f78=46, 5\"),\n    @(\"409\u00f75=81, 4\", \"471\u00f79=52, 3\"),\n    @(\"498\u00f75=99, 3\", \"414\u00f75=82, 4\"),\n    @(\"708\u00f79=78, 6\", \"795\u00f77=113, 4\"),\n    @(\"904\u00f76=150, 4\", \"939\u00f72=469, 1\"),\n    @(\"102\u00f77=14, 4\", \"487\u00f78=60, 7\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
